# update on 20210731 画中人
# Replace curly double quotes (“ ”) with straight single quotes (')
# around quoted phrases in the English (en_US, column C) dialogue text,
# leaving curly apostrophes (’) used for contractions/possessives untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C20").Value2 = "[name=`"Amiya`"]  Our enemy is not only 'Oripathy.' There is no safety for us standing back from this war anymore.`n"
$ws.Range("C34").Value2 = "[name=`"Amiya`"]  ——'The Infected caused this war.'`n"
$ws.Range("C170").Value2 = "[name=`"Firewatch`"]  Those are the 'ordinary Ursus civilians' I know. Had they not supported the war, none of these things would’ve happened.`n"
$ws.Range("C176").Value2 = "[name=`"Firewatch`"]  Who can I trust? How am I supposed to I believe in the benevolence of Ursus’s 'honest folk?'`n"
$ws.Range("C213").Value2 = "[name=`"Dobermann`"]  Hah... We got the wuss who flies the 'Bad Guy' and now this joker.`n"
$ws.Range("C217").Value2 = "[name=`"Pilot`"]  Rest assured, you’re riding the 'Good Boy,' and he’s a lot noisier than the 'Bad Guy!' It’ll be all eyes on us for a solid half hour. Just watch me, Instructor!`n"
